# Apply cryptos list update (prices/volumes refreshed, ImmutableX/NEARProtocol rows swapped)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Donor cell with the default (no explicit style) formatting used to strip the
# text-number-format flag off cells whose new value looks like a plain number,
# so we keep them as Text without altering their style index.
$styleDonor = $ws.Range("B2")

$ws.Range("D2").Value = "65.379.83"
$ws.Range("D3").Value = "3.200.93"
$ws.Range("E3").Value = "  +1.83%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.32"
$styleDonor.Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("E5").Value = "  +1.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.70"
$styleDonor.Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("E6").Value = "  +5.62%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.198.38"
$ws.Range("E8").Value = "  +1.80%  "
$ws.Range("E9").Value = "  +2.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.167"
$styleDonor.Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("E10").Value = "  +3.89%  "
$ws.Range("E11").Value = "  +5.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.472"
$styleDonor.Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("E12").Value = "  +2.91%  "
$ws.Range("E13").Value = "  +3.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.39"
$styleDonor.Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("D15").Value = "3.732.28"
$ws.Range("E15").Value = "  +1.88%  "
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.42"
$styleDonor.Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("E17").Value = "  +4.44%  "
$ws.Range("D18").Value = "65.080.96"
$ws.Range("E18").Value = "  +2.33%  "
$ws.Range("D19").Value = "3.200.09"
$ws.Range("E19").Value = "  +1.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "485.11"
$styleDonor.Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").Value = "  +4.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.13"
$styleDonor.Copy()
$ws.Range("D21").PasteSpecial(-4122)
$ws.Range("E21").Value = "  +6.01%  "
$ws.Range("E22").Value = "  +6.12%  "
$ws.Range("E23").Value = "  +6.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.92"
$styleDonor.Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("E24").Value = "  +7.09%  "
$ws.Range("E25").Value = "  +11.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "83.66"
$styleDonor.Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").Value = "  +2.82%  "
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.86"
$styleDonor.Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = "  +7.32%  "
$ws.Range("E29").Value = "  +3.66%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.51"
$styleDonor.Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = "  +7.61%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.28"
$styleDonor.Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").Value = "  +2.90%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("E33").Value = "  +9.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.61"
$styleDonor.Copy()
$ws.Range("D34").PasteSpecial(-4122)
$ws.Range("E34").Value = "  +6.03%  "
$ws.Range("D35").Value = "0.0₃0896"
$ws.Range("E35").Value = "  +4.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.58"
$styleDonor.Copy()
$ws.Range("D36").PasteSpecial(-4122)
$ws.Range("E36").Value = "  +6.28%  "
$ws.Range("E37").Value = "  +4.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.37"
$styleDonor.Copy()
$ws.Range("D38").PasteSpecial(-4122)
$ws.Range("E38").Value = "  +5.85%  "
$ws.Range("E39").Value = "  +3.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "475.34"
$styleDonor.Copy()
$ws.Range("D40").PasteSpecial(-4122)
$ws.Range("E40").Value = "  +7.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.48"
$styleDonor.Copy()
$ws.Range("D41").PasteSpecial(-4122)
$ws.Range("E41").Value = "  +7.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "51.57"
$styleDonor.Copy()
$ws.Range("D42").PasteSpecial(-4122)
$ws.Range("E42").Value = "  +1.43%  "
$ws.Range("E43").Value = "  +8.31%  "
$ws.Range("E44").Value = "  +3.44%  "
$ws.Range("D45").Value = "2.961.73"
$ws.Range("E45").Value = "  +1.57%  "
$ws.Range("E46").Value = "  +3.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.65"
$styleDonor.Copy()
$ws.Range("D47").PasteSpecial(-4122)
$ws.Range("E47").Value = "  +4.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.73"
$styleDonor.Copy()
$ws.Range("D48").PasteSpecial(-4122)
$ws.Range("E48").Value = "  +4.67%  "
$ws.Range("E49").Value = "  +7.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.69"
$styleDonor.Copy()
$ws.Range("D50").PasteSpecial(-4122)
$ws.Range("E50").Value = "  +5.06%  "

$excel.CutCopyMode = 0
Write-Host "cryptos list updated"
